$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Rule name in the last decision-table row (row 11, column B) changes
# from "R40" to "1". A leading apostrophe forces Excel to keep the entered
# value as text (so it lands in the shared-string table as "1") instead of
# being auto-parsed as the number 1.
$ws.Range("B11").Value = "'1"
